$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1494592.4
$ws.Range("J17").Value = 1540558.9
$ws.Range("L17").Value = 4621676.699999999
$ws.Range("N17").Value = -4622012.699999999
# Row 113
$ws.Range("H113").Value = 4103.15
$ws.Range("I113").Value = 4030.5
$ws.Range("J113").Value = 4175.8
$ws.Range("K113").Value = 4030.5
$ws.Range("L113").Value = 4175.8
$ws.Range("M113").Value = -776.5
$ws.Range("N113").Value = -10683.8
# Row 116
$ws.Range("H116").Value = 4044.125
$ws.Range("I116").Value = 1542.1428
$ws.Range("J116").Value = 5990.1113
$ws.Range("K116").Value = 1542.1428
$ws.Range("L116").Value = 5990.1113
$ws.Range("M116").Value = 1899.8572
$ws.Range("N116").Value = -12874.1113
# Row 129
$ws.Range("H129").Value = 182616.86
$ws.Range("J129").Value = 200844.6
$ws.Range("L129").Value = 602533.8
$ws.Range("N129").Value = -612533.8
# Row 138
$ws.Range("H138").Value = 1925.3662
$ws.Range("I138").Value = 810.75
$ws.Range("J138").Value = 2362.4707
$ws.Range("K138").Value = 2432.25
$ws.Range("L138").Value = 7087.4121
$ws.Range("M138").Value = 2707.75
$ws.Range("N138").Value = -17367.4121

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1277.303
$ws.Range("I2").Value = 1257.8096
$ws.Range("J2").Value = 1311.4166
$ws.Range("K2").Value = 1257.8096
$ws.Range("L2").Value = 1311.4166
$ws.Range("M2").Value = -1144.8096
$ws.Range("N2").Value = -1537.4166
# Row 45
$ws.Range("H45").Value = 2144.0715
$ws.Range("I45").Value = 2368.4375
$ws.Range("K45").Value = 2368.4375
$ws.Range("M45").Value = -1991.4375
# Row 61
$ws.Range("H61").Value = 2304.6897
$ws.Range("I61").Value = 1837.0416
$ws.Range("K61").Value = 1837.0416
$ws.Range("M61").Value = -1625.0416
# Row 63
$ws.Range("H63").Value = 10417999
$ws.Range("I63").Value = 10417999
$ws.Range("K63").Value = 10417999
$ws.Range("M63").Value = -10417313
# Row 66
$ws.Range("H66").Value = 10417999
$ws.Range("I66").Value = 10417999
$ws.Range("K66").Value = 52089995
$ws.Range("M66").Value = -52086563
# Row 116
$ws.Range("H116").Value = 1277.303
$ws.Range("I116").Value = 1257.8096
$ws.Range("J116").Value = 1311.4166
$ws.Range("K116").Value = 1257.8096
$ws.Range("L116").Value = 1311.4166
$ws.Range("M116").Value = 1036.1904
$ws.Range("N116").Value = -5899.4166
# Row 122
$ws.Range("H122").Value = 2514.8333
$ws.Range("I122").Value = 1822.375
$ws.Range("K122").Value = 5467.125
$ws.Range("M122").Value = -3017.125
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
# Row 136
$ws.Range("H136").Value = 2304.6897
$ws.Range("I136").Value = 1837.0416
$ws.Range("K136").Value = 5511.1248
$ws.Range("M136").Value = -2961.1248

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1277.303
$ws.Range("I3").Value = 1257.8096
$ws.Range("J3").Value = 1311.4166
$ws.Range("K3").Value = 1257.8096
$ws.Range("L3").Value = 1311.4166
$ws.Range("M3").Value = -1143.8096
$ws.Range("N3").Value = -1539.4166
# Row 134
$ws.Range("H134").Value = 4310.385
$ws.Range("I134").Value = 4402.8
$ws.Range("K134").Value = 13208.4
$ws.Range("M134").Value = -10673.4

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 639
$ws.Range("I107").Value = 446.57144
$ws.Range("J107").Value = 975.75
$ws.Range("K107").Value = 446.57144
$ws.Range("L107").Value = 975.75
$ws.Range("M107").Value = 1473.42856
$ws.Range("N107").Value = -4815.75

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 5899.1665
$ws.Range("I3").Value = 2198
$ws.Range("K3").Value = 6594
$ws.Range("M3").Value = -6482
# Row 34
$ws.Range("H34").Value = 972.875
$ws.Range("I34").Value = 1000
$ws.Range("J34").Value = 969
$ws.Range("K34").Value = 3000
$ws.Range("L34").Value = 2907
$ws.Range("M34").Value = -2916
$ws.Range("N34").Value = -3075
# Row 60
$ws.Range("H60").Value = 307
$ws.Range("I60").Value = 133.75
$ws.Range("J60").Value = 1000
$ws.Range("K60").Value = 401.25
$ws.Range("L60").Value = 3000
$ws.Range("M60").Value = -150.25
$ws.Range("N60").Value = -3502
# Row 131
$ws.Range("H131").Value = 779.65
$ws.Range("J131").Value = 779.65
$ws.Range("L131").Value = 2338.95
$ws.Range("N131").Value = -12418.95

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3184.8518
$ws.Range("J80").Value = 3328.5881
$ws.Range("L80").Value = 3328.5881
$ws.Range("N80").Value = -5324.5881
# Row 83
$ws.Range("H83").Value = 3184.8518
$ws.Range("J83").Value = 3328.5881
$ws.Range("L83").Value = 16642.9405
$ws.Range("N83").Value = -26626.9405
# Row 99
$ws.Range("H99").Value = 4710
$ws.Range("I99").Value = 4710
$ws.Range("K99").Value = 4710
$ws.Range("M99").Value = -2464

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2804.4375
$ws.Range("I7").Value = 3676
$ws.Range("K7").Value = 3676
$ws.Range("M7").Value = -3564
# Row 22
$ws.Range("H22").Value = 2709.7856
$ws.Range("I22").Value = 2228.0833
$ws.Range("J22").Value = 5600
$ws.Range("K22").Value = 2228.0833
$ws.Range("L22").Value = 5600
$ws.Range("M22").Value = -1933.0833
$ws.Range("N22").Value = -6190
# Row 27
$ws.Range("H27").Value = 2709.7856
$ws.Range("I27").Value = 2228.0833
$ws.Range("J27").Value = 5600
$ws.Range("K27").Value = 2228.0833
$ws.Range("L27").Value = 5600
$ws.Range("M27").Value = -2121.0833
$ws.Range("N27").Value = -5814
# Row 46
$ws.Range("H46").Value = 744.61536
$ws.Range("I46").Value = 661.375
$ws.Range("J46").Value = 877.8
$ws.Range("K46").Value = 661.375
$ws.Range("L46").Value = 877.8
$ws.Range("M46").Value = -473.375
$ws.Range("N46").Value = -1253.8
# Row 61
$ws.Range("H61").Value = 3940.2693
$ws.Range("I61").Value = 2132.7
$ws.Range("J61").Value = 9965.5
$ws.Range("K61").Value = 2132.7
$ws.Range("L61").Value = 9965.5
$ws.Range("M61").Value = -1930.7
$ws.Range("N61").Value = -10369.5
# Row 113
$ws.Range("H113").Value = 3940.2693
$ws.Range("I113").Value = 2132.7
$ws.Range("J113").Value = 9965.5
$ws.Range("K113").Value = 2132.7
$ws.Range("L113").Value = 9965.5
$ws.Range("M113").Value = 37.30000000000018
$ws.Range("N113").Value = -14305.5
# Row 126
$ws.Range("H126").Value = 2804.4375
$ws.Range("I126").Value = 3676
$ws.Range("K126").Value = 11028
$ws.Range("M126").Value = -8558
